# Clean up the "Authors" column (column E) for rows 2-13 on the active sheet.
# The source data packed author records with a run of whitespace after each
# trailing comma; the cleanup widens every such run (2+ spaces) by one extra
# space, leaving any single in-name spaces (e.g. "Imen Ben") untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 13; $r++) {
    $addr = "E" + $r
    $cell = $ws.Range($addr)
    $old = $cell.Value()
    $new = $old -replace '  +', ' $0'
    $cell.Value = $new
}
